# Add direction and asset name into entry point and add vectors list
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) EntryPoint sheet: insert "Direction" (after ID) and "Asset Name" (after
#    Description) columns into the existing table, re-using the table so its
#    id/name ("Table2") is preserved.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("EntryPoint")
$lo = $ws.ListObjects.Item(1)

# Capture the existing row-2 values before we shuffle columns around.
$oldId    = $ws.Range("A2").Value()
$oldDesc  = $ws.Range("B2").Value()
$oldTrust = $ws.Range("C2").Value()
$oldMicro = $ws.Range("D2").Value()

# Column B currently carries the "wrap text" style (it is the Description
# column); strip it before moving things so it does not leak onto the new
# "Direction" column.
$ws.Range("B1:B2").Style = "Normal"

# Grow the table to make room for the two new columns (they land at the end
# for now; we will re-point the header text which re-labels ListColumns).
$lo.Resize($ws.Range("A1:F2"))

# Final header order: ID | Direction | Description | Asset Name | Trust Level | Microservice
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Direction"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Asset Name"
$ws.Range("E1").Value = "Trust Level"
$ws.Range("F1").Value = "Microservice"

# Final data row in the same order. "Asset Name" has no value yet.
$ws.Range("A2").Value = $oldId
$ws.Range("B2").Value = "Exit"
$ws.Range("C2").Value = $oldDesc
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = $oldTrust
$ws.Range("F2").Value = $oldMicro

# Re-apply "wrap text" styling to the Description column (now column C).
$ws.Range("C1:C2").WrapText = $true

# Size the new columns reasonably.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).AutoFit()

# ---------------------------------------------------------------------------
# 2) Add a new "Vectors" worksheet with a one-column table listing the
#    available attack vectors.
# ---------------------------------------------------------------------------
$newws = $wb.Worksheets.Add()
$newws.Name = "Vectors"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newws.Move($null, $lastSheet)

# The Move() call invalidates the earlier object reference in this engine,
# so re-acquire the sheet from the workbook before touching it again.
$vectors = $wb.Worksheets.Item("Vectors")

$vectors.Range("A1").Value = "Name"
$vectors.Range("A2").Value = "Attack vector"

$vlo = $vectors.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $vectors.Range("A1:A2"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$vlo.Name = "Table6"
$vlo.TableStyle = "TableStyleMedium23"

$vectors.Columns.Item(1).AutoFit()
